# Added case for alerts for client users
# Appends new submission-time rows to each of the 4 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Submit orders": add rows 70-73 ---
$ws1 = $wb.Worksheets.Item("Submit orders")

$ws1.Cells.Item(70, 1).Value = "09.28.2022 10:27 (Kyiv+Israel) 07:27 (UTC) 16:27 (Japan) 12:57 (India)"
$ws1.Cells.Item(70, 2).Value = 1.128
$ws1.Cells.Item(70, 3).Value = -0.3149999999999999
$ws1.Cells.Item(70, 4).Value = "***"
$ws1.Cells.Item(70, 5).Value = "***"

$ws1.Cells.Item(71, 1).Value = "09.28.2022 23:21 (Kyiv+Israel) 20:21 (UTC) 05:21 (Japan) 01:51 (India)"
$ws1.Cells.Item(71, 2).Value = 0.9370000000000001
$ws1.Cells.Item(71, 3).Value = -0.1240000000000001
$ws1.Cells.Item(71, 4).Value = "***"
$ws1.Cells.Item(71, 5).Value = "***"

$ws1.Cells.Item(72, 1).Value = "09.29.2022 10:43 (Kyiv+Israel) 07:43 (UTC) 16:43 (Japan) 13:13 (India)"
$ws1.Cells.Item(72, 2).Value = 1.651
$ws1.Cells.Item(72, 3).Value = -0.8380000000000001
$ws1.Cells.Item(72, 4).Value = "***"
$ws1.Cells.Item(72, 5).Value = "***"

$ws1.Cells.Item(73, 1).Value = "09.29.2022 10:44 (Kyiv+Israel) 07:44 (UTC) 16:44 (Japan) 13:14 (India)"
$ws1.Cells.Item(73, 2).Value = "***"
$ws1.Cells.Item(73, 3).Value = "***"
$ws1.Cells.Item(73, 4).Value = 1.821
$ws1.Cells.Item(73, 5).Value = -0.6830000000000001

# --- Sheet "Submit internet survey": add row 70 ---
$ws2 = $wb.Worksheets.Item("Submit internet survey")

$ws2.Cells.Item(70, 1).Value = "09.29.2022 10:47 (Kyiv+Israel) 07:47 (UTC) 16:47 (Japan) 13:17 (India)"
$ws2.Cells.Item(70, 2).Value = "***"
$ws2.Cells.Item(70, 3).Value = "***"
$ws2.Cells.Item(70, 4).Value = 1.258
$ws2.Cells.Item(70, 5).Value = -0.493

# --- Sheet "Submit a phone survey": add row 63 ---
$ws3 = $wb.Worksheets.Item("Submit a phone survey")

$ws3.Cells.Item(63, 1).Value = "09.29.2022 10:49 (Kyiv+Israel) 07:49 (UTC) 16:49 (Japan) 13:19 (India)"
$ws3.Cells.Item(63, 2).Value = "***"
$ws3.Cells.Item(63, 3).Value = "***"
$ws3.Cells.Item(63, 4).Value = 2.008
$ws3.Cells.Item(63, 5).Value = -0.204

# --- Sheet "Checkertificate": add row 65 ---
$ws4 = $wb.Worksheets.Item("Checkertificate")

$ws4.Cells.Item(65, 1).Value = "09.29.2022 10:52 (Kyiv+Israel) 07:52 (UTC) 16:52 (Japan) 13:22 (India)"
$ws4.Cells.Item(65, 2).Value = "***"
$ws4.Cells.Item(65, 3).Value = "***"
$ws4.Cells.Item(65, 4).Value = 2.064
$ws4.Cells.Item(65, 5).Value = -1.142
